$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 <- source row 29
$ws.Range("A8").Value = 111756140
$ws.Range("B8").Value = 89405
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 1202
$ws.Range("F8").Value = "Ullticka"
$ws.Range("G8").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H8").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q8").Value = 453820.6239011836
$ws.Range("R8").Value = 7074037.242731699

# Row 9 <- source row 25
$ws.Range("A9").Value = 111756147
$ws.Range("B9").Value = 89425
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 5442
$ws.Range("F9").Value = "Tallticka"
$ws.Range("G9").Value = "Porodaedalea pini"
$ws.Range("H9").Value = "(Brot.) Murrill"
$ws.Range("Q9").Value = 453989.3915585176
$ws.Range("R9").Value = 7073710.21875874

# Row 10 <- source row 8
$ws.Range("A10").Value = 111756148
$ws.Range("B10").Value = 96266
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 223591
$ws.Range("F10").Value = "Skogsnycklar"
$ws.Range("G10").Value = "Dactylorhiza maculata subsp. fuchsii"
$ws.Range("H10").Value = "(Druce) Hyl."
$ws.Range("Q10").Value = 453747.0542679164
$ws.Range("R10").Value = 7073851.289854143

# Row 11 <- source row 33
$ws.Range("A11").Value = 111756167
$ws.Range("B11").Value = 77515
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = "Garnlav"
$ws.Range("G11").Value = "Alectoria sarmentosa"
$ws.Range("H11").Value = "(Ach.) Ach."
$ws.Range("Q11").Value = 454002.5104495964
$ws.Range("R11").Value = 7073638.391199326

# Row 12 <- source row 34
$ws.Range("A12").Value = 111756155
$ws.Range("B12").Value = 89423
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 5432
$ws.Range("F12").Value = "Granticka"
$ws.Range("G12").Value = "Porodaedalea chrysoloma"
$ws.Range("H12").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q12").Value = 453863.4009631127
$ws.Range("R12").Value = 7073965.428905412

# Row 13 <- source row 22
$ws.Range("A13").Value = 111756168
$ws.Range("B13").Value = 77515
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 6425
$ws.Range("F13").Value = "Garnlav"
$ws.Range("G13").Value = "Alectoria sarmentosa"
$ws.Range("H13").Value = "(Ach.) Ach."
$ws.Range("Q13").Value = 453958.9423245317
$ws.Range("R13").Value = 7073596.134472342

# Row 14 <- source row 13
$ws.Range("A14").Value = 111756170
$ws.Range("B14").Value = 96265
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 219790
$ws.Range("F14").Value = "Fläcknycklar"
$ws.Range("G14").Value = "Dactylorhiza maculata"
$ws.Range("H14").Value = "(L.) Soó"
$ws.Range("Q14").Value = 453738.5427278728
$ws.Range("R14").Value = 7073724.066700204

# Row 16 <- source row 26
$ws.Range("A16").Value = 111756154
$ws.Range("B16").Value = 96674
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 219880
$ws.Range("F16").Value = "Kransrams"
$ws.Range("G16").Value = "Polygonatum verticillatum"
$ws.Range("H16").Value = "(L.) All."
$ws.Range("Q16").Value = 453614.9183513908
$ws.Range("R16").Value = 7074108.35826167

# Row 17 <- source row 10
$ws.Range("A17").Value = 111756163
$ws.Range("B17").Value = 77515
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 6425
$ws.Range("F17").Value = "Garnlav"
$ws.Range("G17").Value = "Alectoria sarmentosa"
$ws.Range("H17").Value = "(Ach.) Ach."
$ws.Range("Q17").Value = 453955.6479769219
$ws.Range("R17").Value = 7073945.9492877

# Row 19 <- source row 17
$ws.Range("A19").Value = 111756161
$ws.Range("B19").Value = 77515
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 6425
$ws.Range("F19").Value = "Garnlav"
$ws.Range("G19").Value = "Alectoria sarmentosa"
$ws.Range("H19").Value = "(Ach.) Ach."
$ws.Range("Q19").Value = 453723.2573215028
$ws.Range("R19").Value = 7074069.623294062

# Row 20 <- source row 19
$ws.Range("A20").Value = 111756166
$ws.Range("B20").Value = 77515
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 6425
$ws.Range("F20").Value = "Garnlav"
$ws.Range("G20").Value = "Alectoria sarmentosa"
$ws.Range("H20").Value = "(Ach.) Ach."
$ws.Range("Q20").Value = 453981.6720900657
$ws.Range("R20").Value = 7073697.065866594

# Row 21 <- source row 24
$ws.Range("A21").Value = 111756143
$ws.Range("B21").Value = 90087
$ws.Range("D21").Value = "LC"
$ws.Range("E21").Value = 3298
$ws.Range("F21").Value = "Trådticka"
$ws.Range("G21").Value = "Climacocystis borealis"
$ws.Range("H21").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q21").Value = 453950.9091414157
$ws.Range("R21").Value = 7073591.829928016

# Row 22 <- source row 20
$ws.Range("A22").Value = 111756139
$ws.Range("B22").Value = 89405
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 1202
$ws.Range("F22").Value = "Ullticka"
$ws.Range("G22").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H22").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q22").Value = 453692.6056797595
$ws.Range("R22").Value = 7074032.491935454

# Row 23 <- source row 28
$ws.Range("A23").Value = 111756159
$ws.Range("B23").Value = 89423
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 5432
$ws.Range("F23").Value = "Granticka"
$ws.Range("G23").Value = "Porodaedalea chrysoloma"
$ws.Range("H23").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q23").Value = 453621.4908246148
$ws.Range("R23").Value = 7073983.576241231

# Row 24 <- source row 30
$ws.Range("A24").Value = 111756150
$ws.Range("B24").Value = 95532
$ws.Range("D24").Value = "LC"
$ws.Range("E24").Value = 221945
$ws.Range("F24").Value = "Revlummer"
$ws.Range("G24").Value = "Lycopodium annotinum"
$ws.Range("H24").Value = "L."
$ws.Range("Q24").Value = 453976.2702886119
$ws.Range("R24").Value = 7073812.112971266

# Row 25 <- source row 31
$ws.Range("A25").Value = 111756172
$ws.Range("B25").Value = 85715
$ws.Range("D25").Value = "NT"
$ws.Range("E25").Value = 510
$ws.Range("F25").Value = "Doftskinn"
$ws.Range("G25").Value = "Cystostereum murrayi"
$ws.Range("H25").Value = "(Berk. & M.A. Curtis.) Pouzar"
$ws.Range("Q25").Value = 453938.5789576455
$ws.Range("R25").Value = 7073959.46382203

# Row 26 <- source row 27
$ws.Range("A26").Value = 111756164
$ws.Range("B26").Value = 77515
$ws.Range("D26").Value = "NT"
$ws.Range("E26").Value = 6425
$ws.Range("F26").Value = "Garnlav"
$ws.Range("G26").Value = "Alectoria sarmentosa"
$ws.Range("H26").Value = "(Ach.) Ach."
$ws.Range("Q26").Value = 453971.0747186596
$ws.Range("R26").Value = 7073820.148138274

# Row 27 <- source row 32
$ws.Range("A27").Value = 111756169
$ws.Range("B27").Value = 77515
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = "Garnlav"
$ws.Range("G27").Value = "Alectoria sarmentosa"
$ws.Range("H27").Value = "(Ach.) Ach."
$ws.Range("Q27").Value = 453910.2023238647
$ws.Range("R27").Value = 7073654.334338664

# Row 28 <- source row 16
$ws.Range("A28").Value = 111756160
$ws.Range("B28").Value = 77515
$ws.Range("D28").Value = "NT"
$ws.Range("E28").Value = 6425
$ws.Range("F28").Value = "Garnlav"
$ws.Range("G28").Value = "Alectoria sarmentosa"
$ws.Range("H28").Value = "(Ach.) Ach."
$ws.Range("Q28").Value = 453815.5156181521
$ws.Range("R28").Value = 7073870.182023689

# Row 29 <- source row 9
$ws.Range("A29").Value = 111756151
$ws.Range("B29").Value = 95532
$ws.Range("D29").Value = "LC"
$ws.Range("E29").Value = 221945
$ws.Range("F29").Value = "Revlummer"
$ws.Range("G29").Value = "Lycopodium annotinum"
$ws.Range("H29").Value = "L."
$ws.Range("Q29").Value = 453609.4901279925
$ws.Range("R29").Value = 7074130.545069677

# Row 30 <- source row 23
$ws.Range("A30").Value = 111756165
$ws.Range("B30").Value = 77515
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 6425
$ws.Range("F30").Value = "Garnlav"
$ws.Range("G30").Value = "Alectoria sarmentosa"
$ws.Range("H30").Value = "(Ach.) Ach."
$ws.Range("Q30").Value = 453984.2379404157
$ws.Range("R30").Value = 7073751.417626478

# Row 31 <- source row 35
$ws.Range("A31").Value = 111756157
$ws.Range("B31").Value = 89423
$ws.Range("D31").Value = "NT"
$ws.Range("E31").Value = 5432
$ws.Range("F31").Value = "Granticka"
$ws.Range("G31").Value = "Porodaedalea chrysoloma"
$ws.Range("H31").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q31").Value = 453981.5111392652
$ws.Range("R31").Value = 7073807.172376178

# Row 32 <- source row 14
$ws.Range("A32").Value = 111756153
$ws.Range("B32").Value = 96674
$ws.Range("D32").Value = "LC"
$ws.Range("E32").Value = 219880
$ws.Range("F32").Value = "Kransrams"
$ws.Range("G32").Value = "Polygonatum verticillatum"
$ws.Range("H32").Value = "(L.) All."
$ws.Range("Q32").Value = 453707.5163784204
$ws.Range("R32").Value = 7073721.869806641

# Row 33 <- source row 21
$ws.Range("A33").Value = 111756162
$ws.Range("B33").Value = 77515
$ws.Range("D33").Value = "NT"
$ws.Range("E33").Value = 6425
$ws.Range("F33").Value = "Garnlav"
$ws.Range("G33").Value = "Alectoria sarmentosa"
$ws.Range("H33").Value = "(Ach.) Ach."
$ws.Range("Q33").Value = 453922.6243923472
$ws.Range("R33").Value = 7073958.370937477

# Row 34 <- source row 11
$ws.Range("A34").Value = 111756158
$ws.Range("B34").Value = 89423
$ws.Range("D34").Value = "NT"
$ws.Range("E34").Value = 5432
$ws.Range("F34").Value = "Granticka"
$ws.Range("G34").Value = "Porodaedalea chrysoloma"
$ws.Range("H34").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q34").Value = 454002.8592168373
$ws.Range("R34").Value = 7073783.424762985

# Row 35 <- source row 12
$ws.Range("A35").Value = 111756142
$ws.Range("B35").Value = 90087
$ws.Range("D35").Value = "LC"
$ws.Range("E35").Value = 3298
$ws.Range("F35").Value = "Trådticka"
$ws.Range("G35").Value = "Climacocystis borealis"
$ws.Range("H35").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q35").Value = 454002.5104495964
$ws.Range("R35").Value = 7073638.391199326
